# Apply cryptocurrency price/volume updates from the Aug 3 2023 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store a numeric-looking value as text,
# matching the original (inline-string) cell type instead of auto-converting
# it to a number.
$q = "`'"

$ws.Range("D2").Value = '29.284.90'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '1.839.85'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("D4").Value = $q + '0.9974'
$ws.Range("D5").Value = $q + '240.94'
$ws.Range("E5").Value = '  -1.34%  '
$ws.Range("D6").Value = $q + '0.6710'
$ws.Range("E6").Value = '  -2.29%  '
$ws.Range("D7").Value = $q + '0.9985'
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("D8").Value = $q + '0.07426'
$ws.Range("E8").Value = '  -1.37%  '
$ws.Range("D9").Value = $q + '0.2964'
$ws.Range("E9").Value = '  -2.20%  '
$ws.Range("D10").Value = $q + '22.99'
$ws.Range("E10").Value = '  -1.28%  '
$ws.Range("D11").Value = $q + '0.07726'
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = $q + '5.037'
$ws.Range("E12").Value = '  -1.12%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = $q + '0.6812'
$ws.Range("E13").Value = '  -0.76%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.767.14'
$ws.Range("E14").Value = '  -4.09%  '
$ws.Range("D15").Value = $q + '86.62'
$ws.Range("E15").Value = '  -2.74%  '
$ws.Range("D16").Value = $q + '6.216'
$ws.Range("E16").Value = '  -1.32%  '
$ws.Range("D17").Value = '29.315.92'
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").Value = $q + '0.000008272'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("D19").Value = $q + '229.90'
$ws.Range("E19").Value = '  -1.85%  '
$ws.Range("D20").Value = $q + '12.59'
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").Value = $q + '0.9985'
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").Value = $q + '7.301'
$ws.Range("E22").Value = '  -3.65%  '
$ws.Range("D23").Value = $q + '0.9991'
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("D24").Value = $q + '160.39'
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D25").Value = $q + '8.741'
$ws.Range("E25").Value = '  -1.41%  '
$ws.Range("D26").Value = $q + '0.1417'
$ws.Range("E26").Value = '  -3.03%  '
$ws.Range("D27").Value = $q + '18.06'
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").Value = $q + '1.511'
$ws.Range("E28").Value = '  -0.85%  '
$ws.Range("D29").Value = $q + '4.219'
$ws.Range("E29").Value = '  -0.25%  '
$ws.Range("D30").Value = $q + '4.097'
$ws.Range("E30").Value = '  -0.68%  '
$ws.Range("D31").Value = $q + '1.197'
$ws.Range("E31").Value = '  -0.56%  '
$ws.Range("D32").Value = $q + '0.05347'
$ws.Range("E32").Value = '  +3.33%  '
$ws.Range("D33").Value = $q + '1.877'
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("D34").Value = $q + '0.7547'
$ws.Range("E34").Value = '  -2.00%  '
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("D36").Value = $q + '2.679'
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").Value = '1.331.75'
$ws.Range("E37").Value = '  +2.73%  '
$ws.Range("D38").Value = $q + '0.01806'
$ws.Range("E38").Value = '  -2.32%  '
$ws.Range("D39").Value = $q + '2.730'
$ws.Range("E39").Value = '  +0.90%  '
$ws.Range("D40").Value = $q + '0.9213'
$ws.Range("E40").Value = '  -2.92%  '
$ws.Range("D41").Value = $q + '5.989'
$ws.Range("E41").Value = '  +4.65%  '
$ws.Range("D42").Value = $q + '0.08294'
$ws.Range("E42").Value = '  +19.97%  '
$ws.Range("D43").Value = $q + '1.000'
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").Value = $q + '103.42'
$ws.Range("E44").Value = '  -1.86%  '
$ws.Range("D45").Value = '1.996.23'
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("E46").Value = '  -0.91%  '
$ws.Range("D47").Value = $q + '0.5171'
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("D48").Value = $q + '64.04'
$ws.Range("E48").Value = '  +1.76%  '
$ws.Range("D49").Value = $q + '1.768'
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").Value = $q + '9.292'
$ws.Range("E50").Value = '  -4.06%  '
$ws.Range("D51").Value = $q + '0.05956'
